$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain plain text (matching the source workbook's
    # inlineStr cells) instead of letting Excel auto-detect numbers/percentages.
    $range.NumberFormat = "@"
    $range.Value = $value
    # Clear the number-format override afterwards so no stray style index is
    # left behind on a cell that originally had no explicit style.
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "302.04"
Set-TextValue $ws.Range("E2") "-0.71%"
Set-TextValue $ws.Range("D3") "31.52"
Set-TextValue $ws.Range("E3") "-3.03%"
Set-TextValue $ws.Range("D4") "5.152"
Set-TextValue $ws.Range("E4") "-2.75%"
Set-TextValue $ws.Range("D5") "0.07415"
Set-TextValue $ws.Range("E5") "-1.05%"
Set-TextValue $ws.Range("D6") "2.146"
Set-TextValue $ws.Range("E6") "41.14%"
Set-TextValue $ws.Range("D7") "7.918"
Set-TextValue $ws.Range("E7") "0.89%"
Set-TextValue $ws.Range("B8") "GateToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D8") "3.770"
Set-TextValue $ws.Range("E8") "-0.75%"
Set-TextValue $ws.Range("B9") "MXToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D9") "0.9267"
Set-TextValue $ws.Range("E9") "1.03%"
Set-TextValue $ws.Range("B10") "WazirX"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1730"
Set-TextValue $ws.Range("E10") "1.70%"
Set-TextValue $ws.Range("B11") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D11") "0.07688"
Set-TextValue $ws.Range("E11") "-2.31%"
Set-TextValue $ws.Range("B12") "MandalaExchangeToken"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.08175"
Set-TextValue $ws.Range("E12") "1.29%"
Set-TextValue $ws.Range("B13") "BitrueCoin"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03035"
Set-TextValue $ws.Range("E13") "0.61%"
Set-TextValue $ws.Range("B14") "BitMartToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09947"
Set-TextValue $ws.Range("E14") "0.32%"
Set-TextValue $ws.Range("B15") "BitForexToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001502"
Set-TextValue $ws.Range("E15") "1.09%"
Set-TextValue $ws.Range("B16") "TigerCash"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.006112"
Set-TextValue $ws.Range("E16") "-3.20%"
Set-TextValue $ws.Range("B17") "LEO"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.466"
Set-TextValue $ws.Range("E17") "-0.33%"
Set-TextValue $ws.Range("D18") "2.232"
Set-TextValue $ws.Range("E18") "0.07%"
Set-TextValue $ws.Range("D20") "0.1348"
Set-TextValue $ws.Range("E20") "0.93%"
Set-TextValue $ws.Range("D21") "4.647"
Set-TextValue $ws.Range("D22") "0.04644"
Set-TextValue $ws.Range("E22") "0.63%"
Set-TextValue $ws.Range("E23") "-2.27%"
Set-TextValue $ws.Range("D24") "0.001223"
Set-TextValue $ws.Range("E24") "0.36%"
Set-TextValue $ws.Range("D25") "0.004493"
Set-TextValue $ws.Range("E25") "0.74%"
Set-TextValue $ws.Range("E26") "-7.07%"
Set-TextValue $ws.Range("D39") "0.01745"
Set-TextValue $ws.Range("E39") "-1.89%"
Set-TextValue $ws.Range("D40") "0.04540"
Set-TextValue $ws.Range("E40") "-0.07%"
Set-TextValue $ws.Range("D41") "0.007094"
Set-TextValue $ws.Range("E41") "-1.25%"
Set-TextValue $ws.Range("D42") "0.1349"
Set-TextValue $ws.Range("E42") "0.08%"
Set-TextValue $ws.Range("D43") "0.002207"
Set-TextValue $ws.Range("E43") "0.09%"
Set-TextValue $ws.Range("D44") "0.01097"
Set-TextValue $ws.Range("E44") "-13.95%"
Set-TextValue $ws.Range("D45") "0.00006274"
Set-TextValue $ws.Range("E45") "1.21%"
Set-TextValue $ws.Range("E46") "-46.05%"
Set-TextValue $ws.Range("D47") "1.928"
Set-TextValue $ws.Range("E47") "2.98%"
